{"js": "// Update the date heading (first paragraph in the document body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst heading = paragraphs.items[0];\nheading.getRange().insertText(\"2025-07-27 Sunday\", Word.InsertLocation.replace);\n\n// Update the five data rows of the single table (0-indexed rows 0, 4, 8, 12, 16;\n// each row has 5 cells, 0-indexed columns 0..4).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = {\n  0: [\"396\u00f76=66, 0\", \"595\u00f79=66, 1\", \"883\u00f74=220, 3\", \"764\u00f76=127, 2\", \"245\u00f77=35, 0\"],\n  4: [\"896\u00f78=112, 0\", \"813\u00f73=271, 0\", \"298\u00f76=49, 4\", \"712\u00f76=118, 4\", \"246\u00f74=61, 2\"],\n  8: [\"285\u00f78=35, 5\", \"790\u00f72=395, 0\", \"622\u00f73=207, 1\", \"969\u00f73=323, 0\", \"730\u00f79=81, 1\"],\n  12: [\"809\u00f75=161, 4\", \"414\u00f77=59, 1\", \"972\u00f79=108, 0\", \"366\u00f74=91, 2\", \"740\u00f79=82, 2\"],\n  16: [\"240\u00f73=80, 0\", \"764\u00f74=191, 0\", \"408\u00f75=81, 3\", \"913\u00f73=304, 1\", \"759\u00f78=94, 7\"],\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const r = parseInt(rowIndex, 10);\n  const values = updates[rowIndex];\n  for (let c = 0; c < values.length; c++) {\n    const cell = table.getCell(r, c);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n    const cellParagraph = cellParagraphs.items[0];\n    cellParagraph.getRange().insertText(values[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading paragraph (first paragraph of the document)\n$d.Paragraphs.Item(1).Range.Text = \"2025-07-27 Sunday\"\n\n# Update the five data rows (rows 1, 5, 9, 13, 17) of the single table.\n# Map of row -> list of new cell values (column order 1..5).\n$t = $d.Tables.Item(1)\n\n$updates = @{\n    1  = @(\"396\u00f76=66, 0\", \"595\u00f79=66, 1\", \"883\u00f74=220, 3\", \"764\u00f76=127, 2\", \"245\u00f77=35, 0\")\n    5  = @(\"896\u00f78=112, 0\", \"813\u00f73=271, 0\", \"298\u00f76=49, 4\", \"712\u00f76=118, 4\", \"246\u00f74=61, 2\")\n    9  = @(\"285\u00f78=35, 5\", \"790\u00f72=395, 0\", \"622\u00f73=207, 1\", \"969\u00f73=323, 0\", \"730\u00f79=81, 1\")\n    13 = @(\"809\u00f75=161, 4\", \"414\u00f77=59, 1\", \"972\u00f79=108, 0\", \"366\u00f74=91, 2\", \"740\u00f79=82, 2\")\n    17 = @(\"240\u00f73=80, 0\", \"764\u00f74=191, 0\", \"408\u00f75=81, 3\", \"913\u00f73=304, 1\", \"759\u00f78=94, 7\")\n}\n\nforeach ($r in $updates.Keys) {\n    $values = $updates[$r]\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $values[$c - 1]\n    }\n}\n"}
